# The table gained two new "line" rows (line7, line8), which pushes the
# existing extr1..extr8 rows down by two rows (they end up on rows 10-17
# instead of 8-15). The sequential index in column A is renumbered to
# stay contiguous (0..15), and the new "line" rows get fresh from_bus /
# to_bus values while the extr rows keep their original from_bus/to_bus
# data, just shifted down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the two brand-new rows (8/9, used by line7/line8) and the two rows
# that extend the sheet beyond its old range (16/17, used by extr7/extr8)
# the same "index column" formatting (bold, centered, bordered) already
# used throughout column A, by copying it from neighboring cells that
# already carry it.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New rows: line7, line8 (in_service = True).
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Shifted-down extr1..extr8 rows (now rows 10-17), renumbering the
# sequential index and restoring each row's own from_bus/to_bus values.
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $false

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $false

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
